{"js": "// Add the new \"bringing in your sprites\" tutorial intro paragraphs and\n// heading, inserted right after the \"Write Up\" title paragraph and before\n// the existing blank paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document starts with: [Title \"Write Up\"], [blank paragraph],\n// [blank Heading 1 paragraph]. Insert the new content before the first\n// blank paragraph so the new paragraphs pick up \"no explicit style\"\n// formatting instead of inheriting the Title style.\nconst anchor = paragraphs.items[1];\n\nconst introPara = anchor.insertParagraph(\n  \"In this tutorial, we will be learning how to bring sprites into our game. Sprites are the faces that will be used by our assets, so that your game elements will be visible inside of the game.\",\n  Word.InsertLocation.before\n);\n\nconst joinPara = anchor.insertParagraph(\n  \"So, if you would like to learn a little more about this, then please join us for our brand-new tutorial this week entitled:\",\n  Word.InsertLocation.before\n);\n\nconst headingPara = anchor.insertParagraph(\n  \"2 Bringing in Your Sprites\",\n  Word.InsertLocation.before\n);\nheadingPara.style = \"Heading 1\";\n\nawait context.sync();\n", "ps1": "# Add the new \"bringing in your sprites\" tutorial intro paragraphs and\n# heading, inserted right after the \"Write Up\" title paragraph and before\n# the existing blank paragraph.\n\n$d = $word.ActiveDocument\n\n# The document starts with: [Title \"Write Up\"], [blank paragraph],\n# [blank Heading 1 paragraph]. Anchor on the blank paragraph (2nd one) so\n# the new text inherits its \"no explicit style\" formatting rather than the\n# Title style of paragraph 1.\n$anchor = $d.Paragraphs.Item(2).Range\n\n# Insert in reverse order since each InsertBefore lands right before the\n# (stationary) anchor point, pushing earlier insertions further down.\n$anchor.InsertBefore(\"2 Bringing in Your Sprites`r\")\n$anchor.InsertBefore(\"So, if you would like to learn a little more about this, then please join us for our brand-new tutorial this week entitled:`r\")\n$anchor.InsertBefore(\"In this tutorial, we will be learning how to bring sprites into our game. Sprites are the faces that will be used by our assets, so that your game elements will be visible inside of the game.`r\")\n\n# The heading paragraph is now the 4th paragraph; give it the Heading 1 style.\n$d.Paragraphs.Item(4).Style = \"Heading 1\"\n"}
